$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows: one before the current row 15 (1-based row 16,
# "Add extended price column...") to hold a duplicated "Check total count..."
# step, and one between the current "Delete Comment column" row and
# "Save file to csv" row to hold a new "Sort by Date, then Item" step.

# Insert a row at row 16 for the extra "Check total count..." step
# (originally row 16 is "Add extended price column and populate with Price * Qty")
$ws.Rows("16:16").Insert()
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = "Check total count against resource file to ensure data is correct"

# Rename the old "Delete Comment column" step (now shifted down to row 20)
$ws.Cells.Item(20, 2).Value = "New data frame with deleted Comment column"

# Insert a new row at row 22 for "Sort by Date, then Item" (after the
# "Check total count..." step that now sits at row 21, before "Save file to csv")
$ws.Rows("22:22").Insert()
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = "Sort by Date, then Item "

# Renumber the remaining Step column values for the final two rows
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(24, 1).Value = 22

# Update the selected cell to match the post-edit active selection
$ws.Range("A25").Select()
